$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 18 (pushes existing rows 18-84 down to 20-86),
# mirroring a new week of price data being prepended to the table.
$ws.Rows.Item(18).Insert()
$ws.Rows.Item(18).Insert()

# New row 18: Cilantro, Primera
$ws.Cells.Item(18,1).Value = 7
$ws.Cells.Item(18,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(18,3).Value = "Ñuble"
$ws.Cells.Item(18,4).Value = 44811
$ws.Cells.Item(18,5).Value = 16
$ws.Cells.Item(18,6).Value = 100112040
$ws.Cells.Item(18,7).Value = "Cilantro"
$ws.Cells.Item(18,8).Value = "Sin especificar"
$ws.Cells.Item(18,9).Value = "Primera"
$ws.Cells.Item(18,10).Value = 200
$ws.Cells.Item(18,11).Value = 700
$ws.Cells.Item(18,12).Value = 800
$ws.Cells.Item(18,13).Value = 750
$ws.Cells.Item(18,14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(18,15).Value = "Provincia de Diguillín"
$ws.Cells.Item(18,16).Value = 750
$ws.Cells.Item(18,17).Value = 1
$ws.Cells.Item(18,18).Value = "Hortaliza"

# New row 19: Cilantro, Segunda
$ws.Cells.Item(19,1).Value = 7
$ws.Cells.Item(19,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(19,3).Value = "Ñuble"
$ws.Cells.Item(19,4).Value = 44811
$ws.Cells.Item(19,5).Value = 16
$ws.Cells.Item(19,6).Value = 100112040
$ws.Cells.Item(19,7).Value = "Cilantro"
$ws.Cells.Item(19,8).Value = "Sin especificar"
$ws.Cells.Item(19,9).Value = "Segunda"
$ws.Cells.Item(19,10).Value = 150
$ws.Cells.Item(19,11).Value = 600
$ws.Cells.Item(19,12).Value = 600
$ws.Cells.Item(19,13).Value = 600
$ws.Cells.Item(19,14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(19,15).Value = "Provincia de Diguillín"
$ws.Cells.Item(19,16).Value = 600
$ws.Cells.Item(19,17).Value = 1
$ws.Cells.Item(19,18).Value = "Hortaliza"
